# Apply weekly fruit/vegetable price update: permute the Fecha/Volumen/Precio/Origen
# data block across data rows 2-117 of the "Puerro" sheet to match the new dataset snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44643
$ws.Range("J2").Value = 160
$ws.Range("K2").Value = 8000
$ws.Range("L2").Value = 9000
$ws.Range("M2").Value = 8500
$ws.Range("P2").Value = 425
$ws.Range("D3").Value = 44286
$ws.Range("J3").Value = 70
$ws.Range("K3").Value = 8000
$ws.Range("M3").Value = 8000
$ws.Range("P3").Value = 400
$ws.Range("D4").Value = 44259
$ws.Range("K4").Value = 8000
$ws.Range("M4").Value = 8000
$ws.Range("P4").Value = 400
$ws.Range("D5").Value = 44215
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 7000
$ws.Range("L5").Value = 7000
$ws.Range("M5").Value = 7000
$ws.Range("P5").Value = 350
$ws.Range("D7").Value = 44349
$ws.Range("J7").Value = 130
$ws.Range("K7").Value = 8000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 8000
$ws.Range("P7").Value = 400
$ws.Range("D8").Value = 44407
$ws.Range("J8").Value = 160
$ws.Range("K8").Value = 7000
$ws.Range("M8").Value = 7500
$ws.Range("P8").Value = 375
$ws.Range("D9").Value = 44873
$ws.Range("J9").Value = 70
$ws.Range("K9").Value = 8000
$ws.Range("L9").Value = 9000
$ws.Range("M9").Value = 8571
$ws.Range("P9").Value = 429
$ws.Range("D10").Value = 44239
$ws.Range("J10").Value = 70
$ws.Range("K10").Value = 8000
$ws.Range("M10").Value = 8000
$ws.Range("P10").Value = 400
$ws.Range("D11").Value = 44664
$ws.Range("J11").Value = 106
$ws.Range("K11").Value = 8000
$ws.Range("M11").Value = 8000
$ws.Range("P11").Value = 400
$ws.Range("D12").Value = 44524
$ws.Range("D13").Value = 44860
$ws.Range("J13").Value = 70
$ws.Range("K13").Value = 9000
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = 9571
$ws.Range("P13").Value = 479
$ws.Range("D14").Value = 44539
$ws.Range("J14").Value = 133
$ws.Range("K14").Value = 6000
$ws.Range("L14").Value = 7000
$ws.Range("M14").Value = 6504
$ws.Range("P14").Value = 325
$ws.Range("D15").Value = 44552
$ws.Range("J15").Value = 106
$ws.Range("K15").Value = 7000
$ws.Range("M15").Value = 7500
$ws.Range("P15").Value = 375
$ws.Range("D16").Value = 44497
$ws.Range("J16").Value = 180
$ws.Range("K16").Value = 6000
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = 6556
$ws.Range("P16").Value = 328
$ws.Range("D17").Value = 44167
$ws.Range("J17").Value = 50
$ws.Range("D18").Value = 44218
$ws.Range("J18").Value = 80
$ws.Range("K18").Value = 6000
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = 6625
$ws.Range("P18").Value = 331
$ws.Range("D19").Value = 44461
$ws.Range("J19").Value = 79
$ws.Range("M19").Value = 7494
$ws.Range("D20").Value = 44720
$ws.Range("J20").Value = 160
$ws.Range("K20").Value = 7000
$ws.Range("M20").Value = 7500
$ws.Range("P20").Value = 375
$ws.Range("D21").Value = 44825
$ws.Range("J21").Value = 70
$ws.Range("K21").Value = 12000
$ws.Range("L21").Value = 12000
$ws.Range("M21").Value = 12000
$ws.Range("P21").Value = 600
$ws.Range("D22").Value = 44748
$ws.Range("J22").Value = 70
$ws.Range("K22").Value = 8000
$ws.Range("M22").Value = 8000
$ws.Range("P22").Value = 400
$ws.Range("D23").Value = 44421
$ws.Range("J23").Value = 180
$ws.Range("D24").Value = 44615
$ws.Range("J24").Value = 79
$ws.Range("K24").Value = 7000
$ws.Range("L24").Value = 7000
$ws.Range("M24").Value = 7000
$ws.Range("P24").Value = 350
$ws.Range("D25").Value = 44229
$ws.Range("J25").Value = 50
$ws.Range("K25").Value = 7000
$ws.Range("M25").Value = 7000
$ws.Range("P25").Value = 350
$ws.Range("D26").Value = 44798
$ws.Range("J26").Value = 52
$ws.Range("K26").Value = 12000
$ws.Range("L26").Value = 13000
$ws.Range("M26").Value = 12500
$ws.Range("P26").Value = 625
$ws.Range("D27").Value = 44519
$ws.Range("K27").Value = 6000
$ws.Range("L27").Value = 7000
$ws.Range("M27").Value = 6500
$ws.Range("P27").Value = 325
$ws.Range("D28").Value = 44244
$ws.Range("J28").Value = 70
$ws.Range("K28").Value = 8000
$ws.Range("M28").Value = 8000
$ws.Range("P28").Value = 400
$ws.Range("D29").Value = 44755
$ws.Range("J29").Value = 160
$ws.Range("K29").Value = 7000
$ws.Range("M29").Value = 7500
$ws.Range("P29").Value = 375
$ws.Range("D30").Value = 44265
$ws.Range("J30").Value = 70
$ws.Range("L30").Value = 8000
$ws.Range("M30").Value = 8000
$ws.Range("P30").Value = 400
$ws.Range("D31").Value = 44707
$ws.Range("K31").Value = 6000
$ws.Range("L31").Value = 7000
$ws.Range("M31").Value = 6571
$ws.Range("P31").Value = 329
$ws.Range("D32").Value = 44398
$ws.Range("J32").Value = 70
$ws.Range("K32").Value = 7500
$ws.Range("L32").Value = 8000
$ws.Range("M32").Value = 7750
$ws.Range("P32").Value = 388
$ws.Range("D33").Value = 44699
$ws.Range("J33").Value = 160
$ws.Range("L33").Value = 8000
$ws.Range("M33").Value = 7500
$ws.Range("P33").Value = 375
$ws.Range("D34").Value = 44224
$ws.Range("J34").Value = 120
$ws.Range("M34").Value = 6667
$ws.Range("P34").Value = 333
$ws.Range("D35").Value = 44189
$ws.Range("J35").Value = 50
$ws.Range("K35").Value = 8000
$ws.Range("L35").Value = 8000
$ws.Range("M35").Value = 8000
$ws.Range("P35").Value = 400
$ws.Range("D36").Value = 44272
$ws.Range("K36").Value = 8000
$ws.Range("M36").Value = 8000
$ws.Range("P36").Value = 400
$ws.Range("D37").Value = 44306
$ws.Range("J37").Value = 160
$ws.Range("D38").Value = 44365
$ws.Range("J38").Value = 180
$ws.Range("D39").Value = 44208
$ws.Range("J39").Value = 50
$ws.Range("K39").Value = 8000
$ws.Range("L39").Value = 8000
$ws.Range("M39").Value = 8000
$ws.Range("P39").Value = 400
$ws.Range("D40").Value = 44510
$ws.Range("J40").Value = 160
$ws.Range("K40").Value = 6000
$ws.Range("L40").Value = 7000
$ws.Range("M40").Value = 6500
$ws.Range("P40").Value = 325
$ws.Range("D41").Value = 44505
$ws.Range("K41").Value = 6000
$ws.Range("L41").Value = 7000
$ws.Range("M41").Value = 6500
$ws.Range("P41").Value = 325
$ws.Range("D42").Value = 44433
$ws.Range("J42").Value = 142
$ws.Range("L42").Value = 8000
$ws.Range("M42").Value = 7500
$ws.Range("P42").Value = 375
$ws.Range("D43").Value = 44419
$ws.Range("D44").Value = 44790
$ws.Range("J44").Value = 160
$ws.Range("K44").Value = 7000
$ws.Range("L44").Value = 7000
$ws.Range("M44").Value = 7000
$ws.Range("P44").Value = 350
$ws.Range("D45").Value = 44266
$ws.Range("J45").Value = 50
$ws.Range("K45").Value = 8000
$ws.Range("M45").Value = 8000
$ws.Range("P45").Value = 400
$ws.Range("D46").Value = 44356
$ws.Range("J46").Value = 160
$ws.Range("K46").Value = 7000
$ws.Range("M46").Value = 7500
$ws.Range("P46").Value = 375
$ws.Range("D48").Value = 44321
$ws.Range("J48").Value = 250
$ws.Range("K48").Value = 7000
$ws.Range("M48").Value = 7000
$ws.Range("P48").Value = 350
$ws.Range("D49").Value = 44358
$ws.Range("K49").Value = 7500
$ws.Range("M49").Value = 7750
$ws.Range("P49").Value = 388
$ws.Range("D50").Value = 44314
$ws.Range("J50").Value = 160
$ws.Range("K50").Value = 8000
$ws.Range("L50").Value = 8000
$ws.Range("M50").Value = 8000
$ws.Range("P50").Value = 400
$ws.Range("D51").Value = 44489
$ws.Range("J51").Value = 160
$ws.Range("K51").Value = 7000
$ws.Range("L51").Value = 8000
$ws.Range("M51").Value = 7500
$ws.Range("P51").Value = 375
$ws.Range("D52").Value = 44273
$ws.Range("J52").Value = 70
$ws.Range("D53").Value = 44855
$ws.Range("J53").Value = 65
$ws.Range("K53").Value = 9000
$ws.Range("L53").Value = 10000
$ws.Range("M53").Value = 9538
$ws.Range("O53").Value = "Provincia de Melipilla"
$ws.Range("P53").Value = 477
$ws.Range("D54").Value = 44855
$ws.Range("J54").Value = 80
$ws.Range("K54").Value = 10000
$ws.Range("L54").Value = 10000
$ws.Range("M54").Value = 10000
$ws.Range("O54").Value = "Provincia de Santiago"
$ws.Range("P54").Value = 500
$ws.Range("D55").Value = 44203
$ws.Range("J55").Value = 50
$ws.Range("K55").Value = 7000
$ws.Range("M55").Value = 7400
$ws.Range("P55").Value = 370
$ws.Range("D56").Value = 44267
$ws.Range("J56").Value = 160
$ws.Range("K56").Value = 8000
$ws.Range("L56").Value = 8000
$ws.Range("M56").Value = 8000
$ws.Range("P56").Value = 400
$ws.Range("D57").Value = 44405
$ws.Range("J57").Value = 160
$ws.Range("D58").Value = 44160
$ws.Range("J58").Value = 50
$ws.Range("M58").Value = 7600
$ws.Range("P58").Value = 380
$ws.Range("D59").Value = 44484
$ws.Range("J59").Value = 160
$ws.Range("L59").Value = 8000
$ws.Range("M59").Value = 7500
$ws.Range("P59").Value = 375
$ws.Range("D60").Value = 44692
$ws.Range("J60").Value = 124
$ws.Range("D61").Value = 44329
$ws.Range("K61").Value = 8000
$ws.Range("M61").Value = 8000
$ws.Range("P61").Value = 400
$ws.Range("D62").Value = 44335
$ws.Range("J62").Value = 250
$ws.Range("K62").Value = 7000
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = 7500
$ws.Range("O62").Value = "Provincia de Chacabuco"
$ws.Range("P62").Value = 375
$ws.Range("D63").Value = 44526
$ws.Range("J63").Value = 170
$ws.Range("K63").Value = 6000
$ws.Range("L63").Value = 7000
$ws.Range("M63").Value = 6500
$ws.Range("O63").Value = "Provincia de Chacabuco"
$ws.Range("P63").Value = 325
$ws.Range("D64").Value = 44302
$ws.Range("J64").Value = 160
$ws.Range("D65").Value = 44195
$ws.Range("K65").Value = 7000
$ws.Range("L65").Value = 7000
$ws.Range("M65").Value = 7000
$ws.Range("P65").Value = 350
$ws.Range("D66").Value = 44370
$ws.Range("K66").Value = 7500
$ws.Range("M66").Value = 7750
$ws.Range("P66").Value = 388
$ws.Range("D67").Value = 44517
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 7000
$ws.Range("M67").Value = 6500
$ws.Range("P67").Value = 325
$ws.Range("D68").Value = 44166
$ws.Range("J68").Value = 50
$ws.Range("D69").Value = 44214
$ws.Range("J69").Value = 50
$ws.Range("D70").Value = 44308
$ws.Range("J70").Value = 160
$ws.Range("K70").Value = 8000
$ws.Range("M70").Value = 8000
$ws.Range("P70").Value = 400
$ws.Range("D71").Value = 44708
$ws.Range("J71").Value = 80
$ws.Range("K71").Value = 7000
$ws.Range("L71").Value = 7000
$ws.Range("M71").Value = 7000
$ws.Range("P71").Value = 350
$ws.Range("D72").Value = 44384
$ws.Range("J72").Value = 160
$ws.Range("M72").Value = 8500
$ws.Range("P72").Value = 425
$ws.Range("D73").Value = 44278
$ws.Range("J73").Value = 70
$ws.Range("K73").Value = 8000
$ws.Range("M73").Value = 8000
$ws.Range("P73").Value = 400
$ws.Range("D74").Value = 44454
$ws.Range("J74").Value = 160
$ws.Range("K74").Value = 7000
$ws.Range("L74").Value = 8000
$ws.Range("M74").Value = 7500
$ws.Range("P74").Value = 375
$ws.Range("D75").Value = 44252
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 160
$ws.Range("K75").Value = 8000
$ws.Range("L75").Value = 8000
$ws.Range("M75").Value = 8000
$ws.Range("O75").Value = "Provincia de Chacabuco"
$ws.Range("P75").Value = 400
$ws.Range("D76").Value = 44475
$ws.Range("D77").Value = 44482
$ws.Range("J77").Value = 160
$ws.Range("K77").Value = 7000
$ws.Range("M77").Value = 7500
$ws.Range("P77").Value = 375
$ws.Range("D78").Value = 44463
$ws.Range("K78").Value = 7500
$ws.Range("M78").Value = 7750
$ws.Range("P78").Value = 388
$ws.Range("D79").Value = 44174
$ws.Range("J79").Value = 70
$ws.Range("K79").Value = 8000
$ws.Range("L79").Value = 8000
$ws.Range("M79").Value = 8000
$ws.Range("P79").Value = 400
$ws.Range("D80").Value = 44503
$ws.Range("J80").Value = 97
$ws.Range("K80").Value = 7000
$ws.Range("M80").Value = 7505
$ws.Range("P80").Value = 375
$ws.Range("D81").Value = 44573
$ws.Range("J81").Value = 106
$ws.Range("D82").Value = 44636
$ws.Range("L82").Value = 7000
$ws.Range("M82").Value = 7000
$ws.Range("P82").Value = 350
$ws.Range("D83").Value = 44301
$ws.Range("J83").Value = 160
$ws.Range("K83").Value = 8000
$ws.Range("L83").Value = 8000
$ws.Range("M83").Value = 8000
$ws.Range("P83").Value = 400
$ws.Range("D84").Value = 44671
$ws.Range("K84").Value = 8000
$ws.Range("M84").Value = 8000
$ws.Range("P84").Value = 400
$ws.Range("D85").Value = 44769
$ws.Range("K85").Value = 7000
$ws.Range("M85").Value = 7500
$ws.Range("P85").Value = 375
$ws.Range("D86").Value = 44232
$ws.Range("J86").Value = 60
$ws.Range("K86").Value = 7000
$ws.Range("L86").Value = 7000
$ws.Range("M86").Value = 7000
$ws.Range("P86").Value = 350
$ws.Range("D87").Value = 44328
$ws.Range("D88").Value = 44442
$ws.Range("J88").Value = 180
$ws.Range("K88").Value = 7000
$ws.Range("L88").Value = 8000
$ws.Range("M88").Value = 7500
$ws.Range("O88").Value = "Provincia de Chacabuco"
$ws.Range("P88").Value = 375
$ws.Range("D89").Value = 44657
$ws.Range("J89").Value = 160
$ws.Range("K89").Value = 7000
$ws.Range("M89").Value = 7500
$ws.Range("P89").Value = 375
$ws.Range("D90").Value = 44498
$ws.Range("K90").Value = 9000
$ws.Range("L90").Value = 10000
$ws.Range("M90").Value = 9615
$ws.Range("O90").Value = "Provincia de Santiago"
$ws.Range("P90").Value = 481
$ws.Range("D91").Value = 44685
$ws.Range("J91").Value = 160
$ws.Range("D93").Value = 44281
$ws.Range("J93").Value = 250
$ws.Range("D94").Value = 44428
$ws.Range("J94").Value = 97
$ws.Range("L94").Value = 9000
$ws.Range("M94").Value = 8505
$ws.Range("P94").Value = 425
$ws.Range("D95").Value = 44363
$ws.Range("J95").Value = 160
$ws.Range("D96").Value = 44344
$ws.Range("J96").Value = 210
$ws.Range("D97").Value = 44426
$ws.Range("J97").Value = 97
$ws.Range("K97").Value = 7000
$ws.Range("L97").Value = 8000
$ws.Range("M97").Value = 7505
$ws.Range("O97").Value = "Provincia de Chacabuco"
$ws.Range("P97").Value = 375
$ws.Range("D98").Value = 44292
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 70
$ws.Range("K98").Value = 8000
$ws.Range("L98").Value = 8000
$ws.Range("M98").Value = 8000
$ws.Range("O98").Value = "Provincia de Chacabuco"
$ws.Range("P98").Value = 400
$ws.Range("D99").Value = 44776
$ws.Range("J99").Value = 160
$ws.Range("M99").Value = 7500
$ws.Range("P99").Value = 375
$ws.Range("D100").Value = 44491
$ws.Range("D101").Value = 44847
$ws.Range("J101").Value = 70
$ws.Range("K101").Value = 11000
$ws.Range("L101").Value = 12000
$ws.Range("M101").Value = 11571
$ws.Range("O101").Value = "Provincia de Melipilla"
$ws.Range("P101").Value = 579
$ws.Range("D102").Value = 44847
$ws.Range("I102").Value = "Segunda"
$ws.Range("J102").Value = 40
$ws.Range("K102").Value = 10000
$ws.Range("M102").Value = 10000
$ws.Range("O102").Value = "Provincia de Melipilla"
$ws.Range("P102").Value = 500
$ws.Range("D103").Value = 44299
$ws.Range("K103").Value = 8000
$ws.Range("L103").Value = 8000
$ws.Range("M103").Value = 8000
$ws.Range("P103").Value = 400
$ws.Range("D104").Value = 44832
$ws.Range("I104").Value = "Segunda"
$ws.Range("J104").Value = 30
$ws.Range("K104").Value = 12000
$ws.Range("L104").Value = 12000
$ws.Range("M104").Value = 12000
$ws.Range("O104").Value = "Provincia de Melipilla"
$ws.Range("P104").Value = 600
$ws.Range("D105").Value = 44162
$ws.Range("J105").Value = 50
$ws.Range("D106").Value = 44468
$ws.Range("J106").Value = 133
$ws.Range("M106").Value = 7504
$ws.Range("D107").Value = 44477
$ws.Range("J107").Value = 160
$ws.Range("K107").Value = 7000
$ws.Range("M107").Value = 7500
$ws.Range("P107").Value = 375
$ws.Range("D108").Value = 44545
$ws.Range("D109").Value = 44295
$ws.Range("J109").Value = 70
$ws.Range("K109").Value = 8000
$ws.Range("L109").Value = 8000
$ws.Range("M109").Value = 8000
$ws.Range("P109").Value = 400
$ws.Range("D110").Value = 44630
$ws.Range("J110").Value = 79
$ws.Range("K110").Value = 9000
$ws.Range("L110").Value = 10000
$ws.Range("M110").Value = 9494
$ws.Range("P110").Value = 475
$ws.Range("D111").Value = 44623
$ws.Range("J111").Value = 106
$ws.Range("K111").Value = 7000
$ws.Range("M111").Value = 7500
$ws.Range("P111").Value = 375
$ws.Range("D112").Value = 44679
$ws.Range("J112").Value = 97
$ws.Range("K112").Value = 8000
$ws.Range("L112").Value = 9000
$ws.Range("M112").Value = 8505
$ws.Range("P112").Value = 425
$ws.Range("D113").Value = 44414
$ws.Range("J113").Value = 180
$ws.Range("K113").Value = 7500
$ws.Range("M113").Value = 7750
$ws.Range("P113").Value = 388
$ws.Range("D114").Value = 44447
$ws.Range("J114").Value = 106
$ws.Range("K114").Value = 7000
$ws.Range("L114").Value = 8000
$ws.Range("M114").Value = 7500
$ws.Range("P114").Value = 375
$ws.Range("D115").Value = 44762
$ws.Range("J115").Value = 160
$ws.Range("K115").Value = 7000
$ws.Range("M115").Value = 7500
$ws.Range("P115").Value = 375
$ws.Range("D116").Value = 44650
$ws.Range("K116").Value = 9000
$ws.Range("L116").Value = 10000
$ws.Range("M116").Value = 9500
$ws.Range("P116").Value = 475
$ws.Range("D117").Value = 44412
$ws.Range("J117").Value = 160
$ws.Range("K117").Value = 7500
$ws.Range("L117").Value = 8000
$ws.Range("M117").Value = 7750
$ws.Range("P117").Value = 388
